# Append CNN training run results (rows 19-25) to the Results_CNN sheet
# as described by the commit: batch normalization, early stopping,
# 48x48 pixels, different models, data augmentation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "2024-1-7 10:27:59"
$ws.Range("B19").Value = 40
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 0.0001
$ws.Range("F19").Value = "CEL"
$ws.Range("G19").Value = 73.40000000000001
$ws.Range("H19").Value = 32
$ws.Range("I19").Value = 0.8885
$ws.Range("J19").Value = 0.51
$ws.Range("K19").Value = 0.5395653385344107
$ws.Range("L19").Value = "FER2013"
$ws.Range("M19").Value = "cpu"
$ws.Range("N19").Value = 4
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = "Alfred"
$ws.Range("Q19").Value = 2935.5
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0.0005
$ws.Range("U19").Value = 0
# Row 20
$ws.Range("A20").Value = "2024-1-7 11:29:54"
$ws.Range("B20").Value = 40
$ws.Range("C20").Value = 64
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = "CEL"
$ws.Range("G20").Value = 68.59999999999999
$ws.Range("H20").Value = 32
$ws.Range("I20").Value = 1.5505
$ws.Range("J20").Value = 1.4351
$ws.Range("K20").Value = 0.4086096405684034
$ws.Range("L20").Value = "FER2013"
$ws.Range("M20").Value = "cpu"
$ws.Range("N20").Value = 4
$ws.Range("O20").Value = 2
$ws.Range("P20").Value = "Alfred"
$ws.Range("Q20").Value = 2745.2
$ws.Range("R20").Value = 0.5
$ws.Range("S20").Value = 0.0001
$ws.Range("T20").Value = "ExponentialLR"
$ws.Range("U20").Value = 0
# Row 21
$ws.Range("A21").Value = "2024-1-7 13:5:31"
$ws.Range("B21").Value = 40
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 0.000048828125
$ws.Range("F21").Value = "CEL"
$ws.Range("G21").Value = 15.8
$ws.Range("H21").Value = 32
$ws.Range("I21").Value = 1.3176
$ws.Range("J21").Value = 1.1116
$ws.Range("K21").Value = 0.4967957648370019
$ws.Range("L21").Value = "FER2013"
$ws.Range("M21").Value = "cpu"
$ws.Range("N21").Value = 4
$ws.Range("O21").Value = 2
$ws.Range("P21").Value = "Alfred"
$ws.Range("Q21").Value = 633.5
$ws.Range("R21").Value = 0.5
$ws.Range("S21").Value = 0.0005
$ws.Range("T21").Value = "AliLR"
$ws.Range("U21").Value = 0
# Row 22
$ws.Range("A22").Value = "2024-1-7 13:29:50"
$ws.Range("B22").Value = 100
$ws.Range("C22").Value = 64
$ws.Range("D22").Value = 0.0000390625
$ws.Range("F22").Value = "CEL"
$ws.Range("G22").Value = 15.7
$ws.Range("H22").Value = 32
$ws.Range("I22").Value = 1.2793
$ws.Range("J22").Value = 0.9649
$ws.Range("K22").Value = 0.5275842853162441
$ws.Range("L22").Value = "FER2013"
$ws.Range("M22").Value = "cpu"
$ws.Range("N22").Value = 4
$ws.Range("O22").Value = 2
$ws.Range("P22").Value = "Alfred"
$ws.Range("Q22").Value = 1272.3
$ws.Range("R22").Value = 0.5
$ws.Range("S22").Value = 0.0005
$ws.Range("T22").Value = "AliLR"
$ws.Range("U22").Value = 0
# Row 23
$ws.Range("A23").Value = "2024-1-7 15:9:56"
$ws.Range("B23").Value = 100
$ws.Range("C23").Value = 64
$ws.Range("D23").Value = 0.000009765625
$ws.Range("F23").Value = "CEL"
$ws.Range("G23").Value = 22.8
$ws.Range("H23").Value = 32
$ws.Range("I23").Value = 1.427
$ws.Range("J23").Value = 0.784
$ws.Range("K23").Value = 0.5125383115073837
$ws.Range("L23").Value = "FER2013"
$ws.Range("M23").Value = "cpu"
$ws.Range("N23").Value = 4
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = "Alfred"
$ws.Range("Q23").Value = 2282.2
$ws.Range("R23").Value = 0.5
$ws.Range("S23").Value = 0.0005
$ws.Range("T23").Value = "AliLR"
$ws.Range("U23").Value = 0
# Row 24
$ws.Range("A24").Value = "2024-1-7 15:39:12"
$ws.Range("B24").Value = 30
$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 0.00015625
$ws.Range("F24").Value = "CEL"
$ws.Range("G24").Value = 67.40000000000001
$ws.Range("H24").Value = 32
$ws.Range("I24").Value = 1.7425
$ws.Range("J24").Value = 1.3977
$ws.Range("K24").Value = 0.3507940930621343
$ws.Range("L24").Value = "FER2013"
$ws.Range("M24").Value = "cpu"
$ws.Range("N24").Value = 4
$ws.Range("O24").Value = 2
$ws.Range("P24").Value = "Alfred"
$ws.Range("Q24").Value = 1281.5
$ws.Range("R24").Value = 0.5
$ws.Range("S24").Value = 0.0005
$ws.Range("T24").Value = "AliLR"
$ws.Range("U24").Value = 0
# Row 25
$ws.Range("A25").Value = "2024-1-7 18:12:38"
$ws.Range("B25").Value = 100
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 0.001
$ws.Range("F25").Value = "CEL"
$ws.Range("G25").Value = 68.5
$ws.Range("H25").Value = 32
$ws.Range("I25").Value = 1.1024
$ws.Range("J25").Value = 0.5744
$ws.Range("K25").Value = 0.6058790749512399
$ws.Range("L25").Value = "FER2013"
$ws.Range("M25").Value = "cpu"
$ws.Range("N25").Value = 4
$ws.Range("O25").Value = 2
$ws.Range("P25").Value = "Alfred"
$ws.Range("Q25").Value = 5340.4
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0.0001
$ws.Range("T25").Value = "None"
$ws.Range("U25").Value = 0
